# edit.ps1 - reproduces the TC002.xlsx edit described by the diff.
#
# Summary of the target state (Sheet1):
#   Row1 (header): Username | PassWord | First Name   -- yellow->green fill, no border
#   Row2 (data)  : DemoCSR  | crmsfa    | M            -- no fill/border
#   Rows 3-5 and the old "Error Message" column content are removed.
#   Column widths / selection / gridlines / window size / page setup updated.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Clear the old content (rows/cols 1:5 / A:C) and rewrite the 2-row table ---
$ws.Cells.Clear()

$ws.Range("A1").Value = "Username"
$ws.Range("B1").Value = "PassWord"
$ws.Range("C1").Value = "First Name"

$ws.Range("A2").Value = "DemoCSR"
$ws.Range("B2").Value = "crmsfa"
$ws.Range("C2").Value = "M"

# --- Header row fill: was yellow (FFFF00), now green (92D050); no border ---
$headerRange = $ws.Range("A1:C1")
$headerRange.Interior.Color = 5296146   # RGB(80,208,146) == BGR 0x92D050 -> matches fgColor FF92D050
$headerRange.Borders.LineStyle = 0      # xlLineStyleNone

# Data row (A2:C2) has no fill and no border in the new style.
$dataRange = $ws.Range("A2:C2")
$dataRange.Interior.Pattern = -4142     # xlPatternNone
$dataRange.Borders.LineStyle = 0        # xlLineStyleNone

# --- Column widths ---
$ws.Columns.Item(1).ColumnWidth = 16.28515625
$ws.Columns.Item(2).ColumnWidth = 9.7109375
$ws.Columns.Item(3).ColumnWidth = 10.5703125

# --- Gridlines off -> on (showGridLines attribute removed from sheetView) ---
$excel.ActiveWindow.DisplayGridlines = $true

# --- Selection moves to D2 ---
$ws.Range("D2").Select()

# --- Page setup: add explicit portrait orientation (adds pageSetup/rId1) ---
$ws.PageSetup.Orientation = 1   # xlPortrait

# --- Workbook window height tweak ---
$excel.ActiveWindow.Height = 601.5
